$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at the top of the "Apio" data
# block (row 70). Every existing record from row 70 downward shifts one
# row down (so the old row 150 becomes row 151), and the dimension grows
# from A1:R150 to A1:R151.

# 1) Push everything from row 70 onward down by one row.
$ws.Rows.Item(70).Insert()

# 2) Seed the brand-new row 70 with the same record that is now sitting
#    in row 71 (i.e. the record that used to live in row 70) ...
$ws.Range("A70:R70").Value2 = $ws.Range("A71:R71").Value2

# 3) ... then apply the two fields that differ for this new weekly entry:
#    the date and the reported volume.
$ws.Range("D70").Value = 44494
$ws.Range("J70").Value = 25
